$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update "Datos actualizados" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 23:52"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 762496
$ws.Range("C4").Value = 23704
$ws.Range("E4").Value = 652062
$ws.Range("G4").Value = 1464
$ws.Range("H4").Value = 40478

# --- España (row 5) ---
$ws.Range("B5").Value = 198674
$ws.Range("C5").Value = 4258
$ws.Range("E5").Value = 100079
$ws.Range("G5").Value = 1195
$ws.Range("H5").Value = 21238

# --- Haiti / Macao swap (rows 162-163) ---
# Haiti moves up to row 162 with updated numbers, Macao moves down to row 163 unchanged
$ws.Range("A162").Value = "Haiti"
$ws.Range("B162").Value = 47
$ws.Range("C162").Value = 3
$ws.Range("D162").Value = 0
$ws.Range("E162").Value = 44
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 0
$ws.Range("H162").Value = 3

$ws.Range("A163").Value = "Macao"
$ws.Range("B163").Value = 45
$ws.Range("C163").Value = 0
$ws.Range("D163").Value = 17
$ws.Range("E163").Value = 28
$ws.Range("F163").Value = 1
$ws.Range("G163").Value = 0
$ws.Range("H163").Value = 0

# --- Gambia / Surinam / Nicaragua reshuffle (rows 200-202) ---
# Gambia moves up to row 200 with updated numbers
$ws.Range("A200").Value = "Gambia"
$ws.Range("B200").Value = 10
$ws.Range("C200").Value = 1
$ws.Range("D200").Value = 2
$ws.Range("E200").Value = 7
$ws.Range("F200").Value = 0
$ws.Range("G200").Value = 0
$ws.Range("H200").Value = 1

# Surinam moves to row 201, unchanged values
$ws.Range("A201").Value = "Surinam"
$ws.Range("B201").Value = 10
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 6
$ws.Range("E201").Value = 3
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 1

# Nicaragua moves to row 202, unchanged values
$ws.Range("A202").Value = "Nicaragua"
$ws.Range("B202").Value = 10
$ws.Range("C202").Value = 1
$ws.Range("D202").Value = 6
$ws.Range("E202").Value = 2
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 2
